# Auto-generated script applying the cryptos.xlsx price/volume update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column cells we are about to touch stay text-typed (the sheet stores
# these as plain strings, e.g. "17.30" or "70.135.45", which must not be reinterpreted as numbers).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated Coin / Link / Price / Volume(1h) values row by row.
$ws.Range("D2").Value = "70.135.45"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "3.484.56"
$ws.Range("E3").Value = "  -2.10%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "615.96"
$ws.Range("E5").Value = "  +2.32%  "
$ws.Range("D6").Value = "168.51"
$ws.Range("E6").Value = "  -2.42%  "
$ws.Range("D7").Value = "3.480.71"
$ws.Range("E7").Value = "  -2.16%  "
$ws.Range("D8").Value = "0.599"
$ws.Range("E8").Value = "  -2.39%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("D11").Value = "7.11"
$ws.Range("E11").Value = "  -4.09%  "
$ws.Range("D12").Value = "0.568"
$ws.Range("E12").Value = "  -3.26%  "
$ws.Range("D13").Value = "44.81"
$ws.Range("E13").Value = "  -3.48%  "
$ws.Range("E14").Value = "  -2.45%  "
$ws.Range("D15").Value = "4.045.64"
$ws.Range("E15").Value = "  -2.26%  "
$ws.Range("D16").Value = "8.25"
$ws.Range("E16").Value = "  -0.94%  "
$ws.Range("D17").Value = "592.20"
$ws.Range("E17").Value = "  -2.67%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.487.06"
$ws.Range("E18").Value = "  -2.11%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "70.243.84"
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("E20").Value = "  +1.07%  "
$ws.Range("D21").Value = "17.30"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").Value = "0.863"
$ws.Range("E22").Value = "  -1.83%  "
$ws.Range("D23").Value = "8.82"
$ws.Range("E23").Value = "  -5.03%  "
$ws.Range("D24").Value = "96.15"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("D25").Value = "15.24"
$ws.Range("E25").Value = "  -2.98%  "
$ws.Range("D26").Value = "3.64"
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "2.49"
$ws.Range("E28").Value = "  -4.41%  "
$ws.Range("D29").Value = "33.18"
$ws.Range("E29").Value = "  -2.41%  "
$ws.Range("D30").Value = "8.72"
$ws.Range("E30").Value = "  -3.85%  "
$ws.Range("D31").Value = "7.93"
$ws.Range("E31").Value = "  -3.60%  "
$ws.Range("D32").Value = "2.84"
$ws.Range("E32").Value = "  -6.78%  "
$ws.Range("E33").Value = "  -2.81%  "
$ws.Range("D34").Value = "6.63"
$ws.Range("E34").Value = "  -5.99%  "
$ws.Range("D35").Value = "572.41"
$ws.Range("E35").Value = "  -19.32%  "
$ws.Range("D36").Value = "0.0489"
$ws.Range("E36").Value = "  +2.08%  "
$ws.Range("D37").Value = "10.72"
$ws.Range("E37").Value = "  -0.36%  "
$ws.Range("D38").Value = "0.0968"
$ws.Range("E38").Value = "  -3.82%  "
$ws.Range("E39").Value = "  +0.29%  "
$ws.Range("D40").Value = "56.29"
$ws.Range("E40").Value = "  -1.19%  "
$ws.Range("E41").Value = "  -0.83%  "
$ws.Range("D42").Value = "3.23"
$ws.Range("E42").Value = "  -10.23%  "
$ws.Range("D43").Value = "3.283.97"
$ws.Range("E43").Value = "  -2.63%  "
$ws.Range("D44").Value = "0.0₃0702"
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("D45").Value = "0.301"
$ws.Range("E45").Value = "  -5.34%  "
$ws.Range("D46").Value = "30.94"
$ws.Range("E46").Value = "  -5.09%  "
$ws.Range("D47").Value = "2.77"
$ws.Range("E47").Value = "  -5.65%  "
$ws.Range("D48").Value = "2.42"
$ws.Range("E48").Value = "  -6.73%  "
$ws.Range("D49").Value = "0.127"
$ws.Range("E49").Value = "  -2.42%  "
$ws.Range("D50").Value = "133.55"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("E51").Value = "  -0.01%  "
